$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of portfolio data (2025-09-30). The date column is stored as
# plain text in this sheet (matching the existing A2:A45 cells), so force
# text interpretation with a leading apostrophe, then strip the resulting
# quote-prefix formatting back to the sheet's default style.
$dateCell = $ws.Range("A46")
$dateCell.Value = "'2025-09-30"
$dateCell.Style = "Normal"

$ws.Range("B46").Value = 55.04999923706055
$ws.Range("C46").Value = 680.2000122070312
$ws.Range("D46").Value = 325.5
